$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.446.66'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.25%  '
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.29%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '484.32'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +4.87%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.91'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.51%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -1.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.998'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.722'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -3.58%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.169'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +8.36%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000354'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +11.98%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '42.44'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -3.58%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.50'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.71%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.545.25'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.91%  '
$ws.Range('B15').NumberFormat = '@'
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').NumberFormat = '@'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.943.99'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.36%  '
$ws.Range('B16').NumberFormat = '@'
$ws.Range('B16').Value = 'Uniswap'
$ws.Range('C16').NumberFormat = '@'
$ws.Range('C16').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.60'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.68%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.39%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '19.69'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -2.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.12'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -3.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '68.496.10'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.30%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '430.98'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.72%  '
$ws.Range('B22').NumberFormat = '@'
$ws.Range('B22').Value = 'InternetComputer(DFINITY)'
$ws.Range('C22').NumberFormat = '@'
$ws.Range('C22').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.48'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -2.89%  '
$ws.Range('B23').NumberFormat = '@'
$ws.Range('B23').Value = 'ImmutableX'
$ws.Range('C23').NumberFormat = '@'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.33'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.39%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '86.88'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.43%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.20'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +11.07%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.33%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.47'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.36%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '38.07'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.91%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.88'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +6.96%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '715.80'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -5.31%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -4.12%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -5.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.83'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +2.96%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0₃0902'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +32.61%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '41.50'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -4.87%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '58.43'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.41%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -7.11%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.47'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.67%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.998'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.17%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.84'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +6.90%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.32%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +11.62%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.99'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.44%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.342'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -4.13%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.140'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.47%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.13%  '
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.24'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.73%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.16%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'LidoDAOToken'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.40'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.62%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '146.98'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.46%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.84'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.83%  '
